$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '29.499.04'
$ws.Range("D2").Style = "Normal"

$ws.Range("E2").Value = '  -2.92%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.997.53'
$ws.Range("D3").Style = "Normal"

$ws.Range("E3").Value = '  -4.61%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.014'
$ws.Range("D4").Style = "Normal"

$ws.Range("E4").Value = '  +1.14%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '329.30'
$ws.Range("D5").Style = "Normal"

$ws.Range("E5").Value = '  -3.89%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.013'
$ws.Range("D6").Style = "Normal"

$ws.Range("E6").Value = '  +1.15%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5007'
$ws.Range("D7").Style = "Normal"

$ws.Range("E7").Value = '  -4.43%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.4222'
$ws.Range("D8").Style = "Normal"

$ws.Range("E8").Value = '  -4.46%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '54.16'
$ws.Range("D9").Style = "Normal"

$ws.Range("E9").Value = '  -0.46%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.08961'
$ws.Range("D10").Style = "Normal"

$ws.Range("E10").Value = '  -4.01%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.114'
$ws.Range("D11").Style = "Normal"

$ws.Range("E11").Value = '  -4.66%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '23.32'
$ws.Range("D12").Style = "Normal"

$ws.Range("E12").Value = '  -5.98%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '2.024.04'
$ws.Range("D13").Style = "Normal"

$ws.Range("E13").Value = '  -4.12%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '7.984'
$ws.Range("D14").Style = "Normal"

$ws.Range("E14").Value = '  -7.18%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.450'
$ws.Range("D15").Style = "Normal"

$ws.Range("E15").Value = '  -6.55%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.013'
$ws.Range("D16").Style = "Normal"

$ws.Range("E16").Value = '  +1.02%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '93.82'
$ws.Range("D17").Style = "Normal"

$ws.Range("E17").Value = '  -7.59%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.00001110'
$ws.Range("D18").Style = "Normal"

$ws.Range("E18").Value = '  -4.20%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06718'
$ws.Range("D19").Style = "Normal"

$ws.Range("E19").Value = '  +0.71%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '19.52'
$ws.Range("D20").Style = "Normal"

$ws.Range("E20").Value = '  -7.85%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '1.015'
$ws.Range("D21").Style = "Normal"

$ws.Range("E21").Value = '  +1.37%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.940'
$ws.Range("D22").Style = "Normal"

$ws.Range("E22").Value = '  -6.22%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '29.529.67'
$ws.Range("D23").Style = "Normal"

$ws.Range("E23").Value = '  -2.86%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '12.00'
$ws.Range("D24").Style = "Normal"

$ws.Range("E24").Value = '  -4.30%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.299'
$ws.Range("D25").Style = "Normal"

$ws.Range("E25").Value = '  -0.62%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '157.50'
$ws.Range("D26").Style = "Normal"

$ws.Range("E26").Value = '  -3.36%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '20.75'
$ws.Range("D27").Style = "Normal"

$ws.Range("E27").Value = '  -5.03%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '6.276'
$ws.Range("D28").Style = "Normal"

$ws.Range("E28").Value = '  -8.11%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.305'
$ws.Range("D29").Style = "Normal"

$ws.Range("E29").Value = '  -7.98%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '127.72'
$ws.Range("D30").Style = "Normal"

$ws.Range("E30").Value = '  -4.20%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.058'
$ws.Range("D31").Style = "Normal"

$ws.Range("E31").Value = '  -6.95%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.09918'
$ws.Range("D32").Style = "Normal"

$ws.Range("E32").Value = '  -5.31%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.559'
$ws.Range("D33").Style = "Normal"

$ws.Range("E33").Value = '  -6.08%  '

$ws.Range("B34").Value = 'Filecoin'

$ws.Range("C34").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.839'
$ws.Range("D34").Style = "Normal"

$ws.Range("E34").Value = '  -6.73%  '

$ws.Range("B35").Value = 'HuobiToken'

$ws.Range("C35").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '3.808'
$ws.Range("D35").Style = "Normal"

$ws.Range("E35").Value = '  -1.31%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.02458'
$ws.Range("D36").Style = "Normal"

$ws.Range("E36").Value = '  -6.60%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '9.267'
$ws.Range("D37").Style = "Normal"

$ws.Range("E37").Value = '  -8.60%  '

$ws.Range("B38").Value = 'Hedera'

$ws.Range("C38").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.06379'
$ws.Range("D38").Style = "Normal"

$ws.Range("E38").Value = '  -6.63%  '

$ws.Range("B39").Value = 'TrustWalletToken'

$ws.Range("C39").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.297'
$ws.Range("D39").Style = "Normal"

$ws.Range("E39").Value = '  -3.15%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.6560'
$ws.Range("D40").Style = "Normal"

$ws.Range("E40").Value = '  -6.11%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '11.61'
$ws.Range("D41").Style = "Normal"

$ws.Range("E41").Value = '  -7.78%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.2043'
$ws.Range("D42").Style = "Normal"

$ws.Range("E42").Value = '  -7.75%  '

$ws.Range("E43").Value = '  +1.26%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.6331'
$ws.Range("D44").Style = "Normal"

$ws.Range("E44").Value = '  -7.31%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '13.53'
$ws.Range("D45").Style = "Normal"

$ws.Range("E45").Value = '  -5.76%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.199'
$ws.Range("D46").Style = "Normal"

$ws.Range("E46").Value = '  -6.21%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.298'
$ws.Range("D47").Style = "Normal"

$ws.Range("E47").Value = '  -6.21%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '3.502'
$ws.Range("D48").Style = "Normal"

$ws.Range("E48").Value = '  -3.83%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.00000000333'
$ws.Range("D49").Style = "Normal"

$ws.Range("E49").Value = '  -5.58%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.06948'
$ws.Range("D50").Style = "Normal"

$ws.Range("E50").Value = '  -4.11%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.129'
$ws.Range("D51").Style = "Normal"

$ws.Range("E51").Value = '  -8.51%  '
